$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B36:F36").Merge()
